$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new data row at row 222 (shifts rows 222:272 down to 223:273)
$ws.Rows.Item(222).Insert()

# Populate the newly inserted row with the new weekly price entry
$ws.Range("A222").Value = 10
$ws.Range("B222").Value = "Vega Modelo de Temuco"
$ws.Range("C222").Value = "La Araucanía"
$ws.Range("D222").Value = 44798
$ws.Range("E222").Value = 9
$ws.Range("F222").Value = 100112043
$ws.Range("G222").Value = "Pepino dulce"
$ws.Range("H222").Value = "Cultivar IV Región"
$ws.Range("I222").Value = "Primera"
$ws.Range("J222").Value = 400
$ws.Range("K222").Value = 18000
$ws.Range("L222").Value = 19000
$ws.Range("M222").Value = 18500
$ws.Range("N222").Value = "`$/bandeja 18 kilos"
$ws.Range("O222").Value = "Provincia de Limarí"
$ws.Range("P222").Value = 1028
$ws.Range("Q222").Value = 18
$ws.Range("R222").Value = "Hortaliza"
